# Slide 19 ("Live demo") has a Content Placeholder with two paragraphs:
#   "(Pending available time)"
#   "https://grey-joyner.netlify.app/"
# Turn the URL text into a real hyperlink pointing at itself, which is
# what PowerPoint does when you select the typed URL and let it
# autoformat / Insert > Link it: it adds an <a:hlinkClick> on the run's
# rPr (via a new external hyperlink relationship) and normalizes the
# paragraph by appending an <a:endParaRPr> that carries the run's
# formatting forward.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(19)
$shape = $s.Shapes.Item(2)

$textRange = $shape.TextFrame.TextRange
$urlParagraph = $textRange.Paragraphs(2)

$hyperlink = $urlParagraph.ActionSettings(1).Hyperlink
$hyperlink.Address = "https://grey-joyner.netlify.app/"
